# Move the "Immagine 40" picture on slide 4 so that it sits before the
# slide title placeholder in the shape (z-order/XML) stack, i.e. send it
# to the back of the stacking order.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$pic = $s.Shapes.Item("Immagine 40")
$pic.ZOrder(1)  # msoSendToBack
